$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-18 00:52:34"
$wsZh.Range("H2").Value = "2016-03-18 00:52:49"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-18 00:52:38"
$wsDe.Range("H2").Value = "2016-03-18 00:52:55"
